$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.627.54'
$ws.Range('E2').Value = '  -5.89%  '
$ws.Range('D3').Value = '1.805.59'
$ws.Range('E3').Value = '  -5.18%  '
$ws.Range('E4').Value = '  +0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '276.35'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  -9.70%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '1.001'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +0.13%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.5060'
$c.Style = 'Normal'
$ws.Range('E7').Value = '  -6.18%  '
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.3520'
$c.Style = 'Normal'
$ws.Range('E8').Value = '  -7.56%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '43.64'
$c.Style = 'Normal'
$ws.Range('E9').Value = '  -5.13%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '0.06635'
$c.Style = 'Normal'
$ws.Range('E10').Value = '  -8.94%  '
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '20.02'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  -9.68%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.8374'
$c.Style = 'Normal'
$ws.Range('E12').Value = '  -7.41%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.07765'
$c.Style = 'Normal'
$ws.Range('E13').Value = '  -5.08%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.815.57'
$ws.Range('E14').Value = '  +54.30%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '5.075'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -5.08%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '87.49'
$c.Style = 'Normal'
$ws.Range('E16').Value = '  -8.59%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +0.04%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '13.93'
$c.Style = 'Normal'
$ws.Range('E18').Value = '  -6.33%  '
$ws.Range('B19').Value = 'Dai'
$ws.Range('C19').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.000007948'
$c.Style = 'Normal'
$ws.Range('E20').Value = '  -8.21%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '25.686.78'
$ws.Range('E21').Value = '  -5.74%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.720'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  -6.49%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '10.02'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  -7.42%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '6.044'
$c.Style = 'Normal'
$ws.Range('E24').Value = '  -7.26%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '142.59'
$c.Style = 'Normal'
$ws.Range('E25').Value = '  -3.95%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '2.108'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  -8.58%  '
$ws.Range('E27').Value = '  -5.66%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '16.91'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -8.08%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '108.41'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  -7.13%  '
$ws.Range('E30').Value = '  -11.20%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '4.216'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  -10.28%  '
$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '0.08798'
$c.Style = 'Normal'
$ws.Range('E32').Value = '  -4.52%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.04787'
$c.Style = 'Normal'
$ws.Range('E33').Value = '  -5.78%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.7236'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  -12.59%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.125'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -7.73%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '2.856'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -4.93%  '
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.9996'
$c.Style = 'Normal'
$ws.Range('E37').Value = '  +0.15%  '
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '3.029'
$c.Style = 'Normal'
$ws.Range('E38').Value = '  -8.83%  '
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.01861'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  -6.94%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.5146'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  -13.39%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.288'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  -15.01%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '0.9614'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  -11.12%  '
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '114.53'
$c.Style = 'Normal'
$ws.Range('E43').Value = '  -1.82%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '6.173'
$c.Style = 'Normal'
$ws.Range('E44').Value = '  -7.38%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '8.030'
$c.Style = 'Normal'
$ws.Range('E45').Value = '  -13.46%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '1.000'
$c.Style = 'Normal'
$ws.Range('E46').Value = '  +0.12%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.4565'
$c.Style = 'Normal'
$ws.Range('E47').Value = '  -10.74%  '
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '0.1381'
$c.Style = 'Normal'
$ws.Range('E48').Value = '  -9.64%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '9.230'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  -9.76%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '35.85'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  -6.28%  '
$ws.Range('E51').Value = '  -9.26%  '
